$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI TPM-derived values for the Col4a5 (ligand) -> Cd93 (receptor) pair.
# Re-running the pipeline with new TPM input changed per-cluster expression stats,
# which ripple into every derived specificity / edge-weight column below.

# Row 2: ECs -> ECs
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.02193933333333333
$ws.Range("H2").Value = 0.065818
$ws.Range("I2").Value = 0.002162808429320595
$ws.Range("J2").Value = 0.002162808429320595
$ws.Range("M2").Value = 83.91225566666667
$ws.Range("N2").Value = 251.736767
$ws.Range("O2").Value = 0.9556261553553385
$ws.Range("P2").Value = 0.9556261553553385
$ws.Range("Q2").Value = 1.840978947822889
$ws.Range("R2").Value = 16.568810530406
$ws.Range("S2").Value = 0.002066836304081759
$ws.Range("T2").Value = 0.002066836304081758

# Row 3: ECs -> FAPs
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.02193933333333333
$ws.Range("H3").Value = 0.065818
$ws.Range("I3").Value = 0.002162808429320595
$ws.Range("J3").Value = 0.002162808429320595
$ws.Range("O3").Value = 0.00439999103960854
$ws.Range("P3").Value = 0.00439999103960854
$ws.Range("Q3").Value = 0.008476422321777779
$ws.Range("R3").Value = 0.07628780089600001
$ws.Range("S3").Value = [double]"9.516337709400439E-06"
$ws.Range("T3").Value = [double]"9.516337709400438E-06"

# Row 4: ECs -> MuSCs
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.02193933333333333
$ws.Range("H4").Value = 0.065818
$ws.Range("I4").Value = 0.002162808429320595
$ws.Range("J4").Value = 0.002162808429320595
$ws.Range("M4").Value = 3.510050666666667
$ws.Range("N4").Value = 10.530152
$ws.Range("O4").Value = 0.03997385360505296
$ws.Range("P4").Value = 0.03997385360505297
$ws.Range("Q4").Value = 0.0770081715928889
$ws.Range("R4").Value = 0.6930735443360001
$ws.Range("S4").Value = [double]"8.6455787529436E-05"
$ws.Range("T4").Value = [double]"8.6455787529436E-05"

# Row 5: FAPs -> ECs
$ws.Range("G5").Value = 7.730541000000001
$ws.Range("I5").Value = 0.7620869323593149
$ws.Range("J5").Value = 0.7620869323593149
$ws.Range("M5").Value = 83.91225566666667
$ws.Range("N5").Value = 251.736767
$ws.Range("O5").Value = 0.9556261553553385
$ws.Range("P5").Value = 0.9556261553553385
$ws.Range("Q5").Value = 648.6871328336491
$ws.Range("R5").Value = 5838.184195502842
$ws.Range("S5").Value = 0.728270205217076
$ws.Range("T5").Value = 0.728270205217076

# Row 6: FAPs -> FAPs
$ws.Range("G6").Value = 7.730541000000001
$ws.Range("I6").Value = 0.7620869323593149
$ws.Range("J6").Value = 0.7620869323593149
$ws.Range("O6").Value = 0.00439999103960854
$ws.Range("P6").Value = 0.00439999103960854
$ws.Range("Q6").Value = 2.986751205984001
$ws.Range("R6").Value = 26.88076085385601
$ws.Range("S6").Value = 0.003353175673783745
$ws.Range("T6").Value = 0.003353175673783745

# Row 7: FAPs -> MuSCs
$ws.Range("G7").Value = 7.730541000000001
$ws.Range("I7").Value = 0.7620869323593149
$ws.Range("J7").Value = 0.7620869323593149
$ws.Range("M7").Value = 3.510050666666667
$ws.Range("N7").Value = 10.530152
$ws.Range("O7").Value = 0.03997385360505296
$ws.Range("P7").Value = 0.03997385360505297
$ws.Range("Q7").Value = 27.13459059074401
$ws.Range("R7").Value = 244.2113153166961
$ws.Range("S7").Value = 0.03046355146845516
$ws.Range("T7").Value = 0.03046355146845516

# Row 8: MuSCs -> ECs
$ws.Range("G8").Value = 2.391429333333333
$ws.Range("H8").Value = 7.174287999999999
$ws.Range("I8").Value = 0.2357502592113645
$ws.Range("J8").Value = 0.2357502592113645
$ws.Range("M8").Value = 83.91225566666667
$ws.Range("N8").Value = 251.736767
$ws.Range("O8").Value = 0.9556261553553385
$ws.Range("P8").Value = 0.9556261553553385
$ws.Range("Q8").Value = 200.6702296274329
$ws.Range("R8").Value = 1806.032066646896
$ws.Range("S8").Value = 0.2252891138341807
$ws.Range("T8").Value = 0.2252891138341807

# Row 9: MuSCs -> FAPs
$ws.Range("G9").Value = 2.391429333333333
$ws.Range("H9").Value = 7.174287999999999
$ws.Range("I9").Value = 0.2357502592113645
$ws.Range("J9").Value = 0.2357502592113645
$ws.Range("O9").Value = 0.00439999103960854
$ws.Range("P9").Value = 0.00439999103960854
$ws.Range("Q9").Value = 0.9239462600817777
$ws.Range("R9").Value = 8.315516340736
$ws.Range("S9").Value = 0.001037299028115395
$ws.Range("T9").Value = 0.001037299028115395

# Row 10: MuSCs -> MuSCs
$ws.Range("G10").Value = 2.391429333333333
$ws.Range("H10").Value = 7.174287999999999
$ws.Range("I10").Value = 0.2357502592113645
$ws.Range("J10").Value = 0.2357502592113645
$ws.Range("M10").Value = 3.510050666666667
$ws.Range("N10").Value = 10.530152
$ws.Range("O10").Value = 0.03997385360505296
$ws.Range("P10").Value = 0.03997385360505297
$ws.Range("Q10").Value = 8.394038125752887
$ws.Range("R10").Value = 75.546343131776
$ws.Range("S10").Value = 0.009423846349068373
$ws.Range("T10").Value = 0.009423846349068375
